$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Name" column values for rows 5 and 6
$ws.Range("A5").Value = "TT03"
$ws.Range("A6").Value = "DP04"
